# Insert a new daily price record for "Ajo" (Feria Lagunitas de Puerto Montt)
# at row 130, pushing the existing rows 130-250 down to 131-251.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 130 (shifts rows 130..250 -> 131..251,
# carrying their formatting/styles along, same as Excel's Insert Row).
$ws.Rows(130).Insert()

# Populate the newly inserted row 130 with the new record.
$ws.Cells.Item(130, 1).Value  = 4
$ws.Cells.Item(130, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(130, 3).Value  = "Los Lagos"
$ws.Cells.Item(130, 4).Value  = 44669
$ws.Cells.Item(130, 5).Value  = 10
$ws.Cells.Item(130, 6).Value  = 100112003
$ws.Cells.Item(130, 7).Value  = "Ajo"
$ws.Cells.Item(130, 8).Value  = "Chino"
$ws.Cells.Item(130, 9).Value  = "Primera"
$ws.Cells.Item(130, 10).Value = 20
$ws.Cells.Item(130, 11).Value = 20000
$ws.Cells.Item(130, 12).Value = 21000
$ws.Cells.Item(130, 13).Value = 20500
$ws.Cells.Item(130, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(130, 15).Value = "China"
$ws.Cells.Item(130, 16).Value = 2050
$ws.Cells.Item(130, 17).Value = 10
$ws.Cells.Item(130, 18).Value = "Hortaliza"
